$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B2").Value = 15.10031222304571
$ws.Range("C2").Value = 11.9368833554878
$ws.Range("D2").Value = 5.131882537702833
$ws.Range("E2").Value = 9.893039654352568
$ws.Range("F2").Value = 30.69742812745123
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 22.91209905447385
$ws.Range("M2").Value = 16.21238624863645
$ws.Range("N2").Value = 17.73480112790808

$ws.Range("B3").Value = 14.49656486284214
$ws.Range("C3").Value = 11.28452866778956
$ws.Range("D3").Value = 5.158064285911068
$ws.Range("E3").Value = 9.806446140188788
$ws.Range("F3").Value = 30.32580373260764
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 22.86487911953112
$ws.Range("M3").Value = 15.92084982569947
$ws.Range("N3").Value = 17.80463104271796

$ws.Range("B4").Value = 14.11715935678368
$ws.Range("C4").Value = 10.86731478749963
$ws.Range("D4").Value = 5.174906245129847
$ws.Range("E4").Value = 9.756137006042151
$ws.Range("F4").Value = 30.10653554502547
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 22.84271331921616
$ws.Range("M4").Value = 15.74425087070024
$ws.Range("N4").Value = 17.8494530844848

$ws.Range("B5").Value = 13.96062383440401
$ws.Range("C5").Value = 10.69331013941727
$ws.Range("D5").Value = 5.181962297801277
$ws.Range("E5").Value = 9.736371651988062
$ws.Range("F5").Value = 30.01952353422046
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 22.83539627480157
$ws.Range("M5").Value = 15.67299155904073
$ws.Range("N5").Value = 17.86820881083082

$ws.Range("B6").Value = 13.93452323559178
$ws.Range("C6").Value = 10.66418272942907
$ws.Range("D6").Value = 5.183145599991938
$ws.Range("E6").Value = 9.733134546498187
$ws.Range("F6").Value = 30.00521950941343
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 22.83428483949116
$ws.Range("M6").Value = 15.66120460603631
$ws.Range("N6").Value = 17.87135283214109

$ws.Range("B7").Value = 14.11505568818162
$ws.Range("C7").Value = 10.86498394410646
$ws.Range("D7").Value = 5.175000624659633
$ws.Range("E7").Value = 9.755867442278859
$ws.Range("F7").Value = 30.10535246115369
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 22.84260769526414
$ws.Range("M7").Value = 15.74328684994622
$ws.Range("N7").Value = 17.84970404399676

$ws.Range("B8").Value = 14.89409574379577
$ws.Range("C8").Value = 11.7155215764403
$ws.Range("D8").Value = 5.140751045699918
$ws.Range("E8").Value = 9.862599363734661
$ws.Range("F8").Value = 30.56750493762174
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 22.89440047680648
$ws.Range("M8").Value = 16.11142606986684
$ws.Range("N8").Value = 17.75847515111597

$ws.Range("B9").Value = 16.3425676286459
$ws.Range("C9").Value = 13.2439666925263
$ws.Range("D9").Value = 5.079661462064228
$ws.Range("E9").Value = 10.09376866516837
$ws.Range("F9").Value = 31.53978443110541
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 23.05010523997703
$ws.Range("M9").Value = 16.84798749576562
$ws.Range("N9").Value = 17.59497098773458

$ws.Range("B10").Value = 17.34669929272451
$ws.Range("C10").Value = 14.27425869691708
$ws.Range("D10").Value = 5.038473088112043
$ws.Range("E10").Value = 10.27574289181346
$ws.Range("F10").Value = 32.28775897434545
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 23.19734239285955
$ws.Range("M10").Value = 17.39219442637149
$ws.Range("N10").Value = 17.48416269468315

$ws.Range("B11").Value = 17.7884680547886
$ws.Range("C11").Value = 14.72170801620016
$ws.Range("D11").Value = 5.020536217572721
$ws.Range("E11").Value = 10.36088360206437
$ws.Range("F11").Value = 32.63387970454265
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 23.27138062324921
$ws.Range("M11").Value = 17.63923393746358
$ws.Range("N11").Value = 17.43576244757029

$ws.Range("B12").Value = 17.95345810135597
$ws.Range("C12").Value = 14.8880211107172
$ws.Range("D12").Value = 5.013859008319539
$ws.Range("E12").Value = 10.39343879107897
$ws.Range("F12").Value = 32.76566636058784
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 23.30042263300957
$ws.Range("M12").Value = 17.73261094910247
$ws.Range("N12").Value = 17.4177221890837

$ws.Range("B13").Value = 17.91802876287642
$ws.Range("C13").Value = 14.8523426483305
$ws.Range("D13").Value = 5.015291944566405
$ws.Range("E13").Value = 10.38641387772311
$ws.Range("F13").Value = 32.73725378515827
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 23.29412337915821
$ws.Range("M13").Value = 17.71250965352261
$ws.Range("N13").Value = 17.42159469009595

$ws.Range("B14").Value = 17.80208867152743
$ws.Range("C14").Value = 14.73545371180735
$ws.Range("D14").Value = 5.019984572403398
$ws.Range("E14").Value = 10.3635557847337
$ws.Range("F14").Value = 32.64470819100374
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 23.2737498409758
$ws.Range("M14").Value = 17.64692009366238
$ws.Range("N14").Value = 17.43427249954306

$ws.Range("B15").Value = 17.7307689330302
$ws.Range("C15").Value = 14.6634468360975
$ws.Range("D15").Value = 5.022873933597618
$ws.Range("E15").Value = 10.34959471301081
$ws.Range("F15").Value = 32.58811118004633
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 23.26140106349893
$ws.Range("M15").Value = 17.60671947622998
$ws.Range("N15").Value = 17.44207549093076

$ws.Range("B16").Value = 17.31751418241527
$ws.Range("C16").Value = 14.24458295866182
$ws.Range("D16").Value = 5.039661406295625
$ws.Range("E16").Value = 10.2702241005561
$ws.Range("F16").Value = 32.26524647311798
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 23.1926450779839
$ws.Range("M16").Value = 17.37603096345278
$ws.Range("N16").Value = 17.4873660867198

$ws.Range("B17").Value = 17.06004470108136
$ws.Range("C17").Value = 13.98212794280221
$ws.Range("D17").Value = 5.0501648948204
$ws.Range("E17").Value = 10.22211954678265
$ws.Range("F17").Value = 32.06859225904245
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 23.1522671986294
$ws.Range("M17").Value = 17.23430727059769
$ws.Range("N17").Value = 17.51566390535942

$ws.Range("B18").Value = 16.91054869217514
$ws.Range("C18").Value = 13.82917521644658
$ws.Range("D18").Value = 5.056281510803138
$ws.Range("E18").Value = 10.19467473560993
$ws.Range("F18").Value = 31.95604133489249
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 23.12970807191559
$ws.Range("M18").Value = 17.15274784033753
$ws.Range("N18").Value = 17.53212901304462

$ws.Range("B19").Value = 16.85969492433977
$ws.Range("C19").Value = 13.77704778134629
$ws.Range("D19").Value = 5.058365420992613
$ws.Range("E19").Value = 10.18542156080332
$ws.Range("F19").Value = 31.91803340894777
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 23.12218444881247
$ws.Range("M19").Value = 17.12512871338252
$ws.Range("N19").Value = 17.53773629107299

$ws.Range("B20").Value = 17.08759938987901
$ws.Range("C20").Value = 14.01027383262416
$ws.Range("D20").Value = 5.049038988221686
$ws.Range("E20").Value = 10.22721739081188
$ws.Range("F20").Value = 32.08946944813473
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 23.1564967161115
$ws.Range("M20").Value = 17.24939923441657
$ws.Range("N20").Value = 17.51263200402581

$ws.Range("B21").Value = 17.83620647049101
$ws.Range("C21").Value = 14.76987217563928
$ws.Range("D21").Value = 5.018603109551486
$ws.Range("E21").Value = 10.37026143482955
$ws.Range("F21").Value = 32.67187257529869
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 23.27970684548149
$ws.Range("M21").Value = 17.66619073768118
$ws.Range("N21").Value = 17.43054091276562

$ws.Range("B22").Value = 18.31201763722913
$ws.Range("C22").Value = 15.24806781961368
$ws.Range("D22").Value = 4.999382413771605
$ws.Range("E22").Value = 10.4655686734462
$ws.Range("F22").Value = 33.05663579176586
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 23.36608584714967
$ws.Range("M22").Value = 17.93754694864518
$ws.Range("N22").Value = 17.37856721350891

$ws.Range("B23").Value = 18.05933963668492
$ws.Range("C23").Value = 14.99453551045115
$ws.Range("D23").Value = 5.009579442744566
$ws.Range("E23").Value = 10.41454333094003
$ws.Range("F23").Value = 32.85094370068454
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 23.31945180583704
$ws.Range("M23").Value = 17.79284498202968
$ws.Range("N23").Value = 17.40615329304031

$ws.Range("B24").Value = 17.07514649355948
$ws.Range("C24").Value = 13.99755549858193
$ws.Range("D24").Value = 5.049547767582985
$ws.Range("E24").Value = 10.2249119959853
$ws.Range("F24").Value = 32.08002927852922
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 23.15458251143727
$ws.Range("M24").Value = 17.2425764012462
$ws.Range("N24").Value = 17.51400211508582

$ws.Range("B25").Value = 15.96051161177802
$ws.Range("C25").Value = 12.84629915081959
$ws.Range("D25").Value = 5.095538561025974
$ws.Range("E25").Value = 10.02900877599982
$ws.Range("F25").Value = 31.27035968232338
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 23.0021962433866
$ws.Range("M25").Value = 16.64777455586022
$ws.Range("N25").Value = 17.63756206291654
